# Update countries & provincias Spain
# Applies the data refresh captured in the commit "Update countries & provincias Spain":
#  - swaps the display order of a couple of country-name pairs whose rows were
#    re-sorted upstream (Groenlandia/Islas Malvinas and Montserrat/Seychelles),
#  - refreshes the numeric columns (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for the countries whose
#    figures moved,
#  - bumps the "Datos actualizados" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 10:24"

# --- Country label swaps (rows kept in place, labels exchanged) -------
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"

$ws.Range("A211").Value = "Seychelles"
$ws.Range("A212").Value = "Montserrat"

# --- Numeric refresh ----------------------------------------------------
# Row 6 - Rusia
$ws.Range("B6").Value = 613994
$ws.Range("C6").Value = 7113
$ws.Range("D6").Value = 375164
$ws.Range("E6").Value = 230225
$ws.Range("G6").Value = 92
$ws.Range("H6").Value = 8605

# Row 7 - India
$ws.Range("B7").Value = 474272
$ws.Range("C7").Value = 1287
$ws.Range("D7").Value = 271934
$ws.Range("E7").Value = 187424
$ws.Range("G7").Value = 7
$ws.Range("H7").Value = 14914

# Row 35 - Singapur
$ws.Range("B35").Value = 42736
$ws.Range("C35").Value = 113
$ws.Range("E35").Value = 6411

# Row 41 - Polonia
$ws.Range("D41").Value = 18654
$ws.Range("E41").Value = 12771

# Row 57 - Moldavia
$ws.Range("D57").Value = 8599
$ws.Range("E57").Value = 5977
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 502

# Row 68 - Chequia
$ws.Range("B68").Value = 10780
$ws.Range("C68").Value = 3
$ws.Range("D68").Value = 7592
$ws.Range("E68").Value = 2844
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 344

# Row 108 - Sri Lanka
$ws.Range("D108").Value = 1602
$ws.Range("E108").Value = 388

# Row 109 - Estonia
$ws.Range("B109").Value = 1984
$ws.Range("C109").Value = 1
$ws.Range("D109").Value = 1790
$ws.Range("E109").Value = 125

# Row 112 - Lituania
$ws.Range("B112").Value = 1806
$ws.Range("C112").Value = 2
$ws.Range("D112").Value = 1494
$ws.Range("E112").Value = 234

# Row 116 - Eslovaquia
$ws.Range("B116").Value = 1630
$ws.Range("C116").Value = 23
$ws.Range("D116").Value = 1452
$ws.Range("E116").Value = 150

# Row 155 - Taiwan
$ws.Range("B155").Value = 447
$ws.Range("C155").Value = 1
$ws.Range("E155").Value = 5

# Row 211 (now Seychelles)
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 212 (now Montserrat)
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
